$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update progress values so tasks 3-6 show as complete (100%)
$ws.Range("B3").Value = 100
$ws.Range("B4").Value = 100
$ws.Range("B5").Value = 100
$ws.Range("B6").Value = 100

# Remove the "배포 준비" task row content (A7 and C7), leaving B7 = 0
$ws.Range("A7").ClearContents()
$ws.Range("C7").ClearContents()

# Update the selected cell as recorded in the saved view
$ws.Range("C10").Select()
